# Updated Masterdata as per 2nd may Data Refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data refresh: update device/center id (column A) for several rows ---
# Rows that held 10002 -> 10003
$rowsFrom10002 = @(3, 23, 43, 63, 83)
foreach ($r in $rowsFrom10002) {
    $ws.Cells.Item($r, 1).Value = 10003
}

# Rows that held 10005 -> 10003
$rowsFrom10005 = @(105, 114, 123, 132, 141)
foreach ($r in $rowsFrom10005) {
    $ws.Cells.Item($r, 1).Value = 10003
}

# --- Restore view to top of sheet and select from row 162 to the last row ---
$ws.Application.Goto($ws.Range("A1"), $true)

$topRow = $ws.Rows.Item(162)
$lastRow = $ws.Rows.Item(1048576)
$selRange = $ws.Range($topRow, $lastRow)
$selRange.Select()
